$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 2) with the lawsuit record.
# The "Filing Date" column holds a plain text date string (matching the
# existing header/data cells which are all stored as text), so it is
# entered with a leading apostrophe to keep Excel from auto-converting it
# into a date serial number.
$ws.Range("A2").Value = "1998D000001"
$ws.Range("B2").Value = "Rojas Nancy"
$ws.Range("C2").Value = "Cochran Oscar"
$ws.Range("D2").Value = "'1998-02-01"
